$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE (column N) for rows 2-5 with the recalculated value
$ws.Range("N2:N5").Value = 85.77505782882612
